$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item(1).Name = "Seating"
$wb.Worksheets.Item(2).Name = "Contact"
$wb.Worksheets.Item(3).Name = "Car"

$wsSeating = $wb.Worksheets.Item("Seating")
$wsContact = $wb.Worksheets.Item("Contact")
$wsCar = $wb.Worksheets.Item("Car")

# --- Contact sheet data ---
$wsContact.Range("A1").Value = "伴娘 "
$wsContact.Range("A1").Characters(3,1).Font.Name = "Times New Roman"
$wsContact.Range("A1").Characters(3,1).Font.Size = 12
$wsContact.Range("B1").Value = "嘉欣"
$wsContact.Range("C1").Value = "9804 6277"
$wsContact.Range("E1").Value = "伴郎 "
$wsContact.Range("E1").Characters(3,1).Font.Name = "Times New Roman"
$wsContact.Range("E1").Characters(3,1).Font.Size = 12
$wsContact.Range("F1").Value = "林祺"
$wsContact.Range("G1").Value = "9802 9700"

$wsContact.Range("A2").Value = "MC/姊妹"
$wsContact.Range("B2").Value = "C 吻"
$wsContact.Range("C2").Value = "9747 1576"
$wsContact.Range("E2").Value = "MC/兄弟"
$wsContact.Range("F2").Value = "Edwin"
$wsContact.Range("G2").Value = "6333 8193"

$wsContact.Range("A3").Value = "姊妹"
$wsContact.Range("B3").Value = "Bus"
$wsContact.Range("C3").Value = "6010 5191"
$wsContact.Range("E3").Value = "兄弟"
$wsContact.Range("F3").Value = "Francis"
$wsContact.Range("G3").Value = "9608 4620"

$wsContact.Range("A4").Value = "姊妹"
$wsContact.Range("B4").Value = "Mouse"
$wsContact.Range("C4").Value = "6229 4669"
$wsContact.Range("E4").Value = "兄弟"
$wsContact.Range("F4").Value = "細 Tony"
$wsContact.Range("G4").Value = "5182 1122"

$wsContact.Range("A5").Value = "姊妹"
$wsContact.Range("B5").Value = "趣影"
$wsContact.Range("C5").Value = "6686 1409"
$wsContact.Range("E5").Value = "兄弟"
$wsContact.Range("F5").Value = "Yellow"
$wsContact.Range("G5").Value = "9226 1217"

$wsContact.Range("A6").Value = "姊妹"
$wsContact.Range("B6").Value = "詩薇"
$wsContact.Range("C6").Value = "6341 3634"
$wsContact.Range("E6").Value = "兄弟"
$wsContact.Range("F6").Value = "Jonathan"
$wsContact.Range("G6").Value = "9529 1899"

$wsContact.Range("A7").Value = "姊妹"
$wsContact.Range("B7").Value = "小如"
$wsContact.Range("C7").Value = "9792 3818"
$wsContact.Range("E7").Value = "兄弟"
$wsContact.Range("F7").Value = "Anthony"
$wsContact.Range("G7").Value = "6409 3285"

$wsContact.Range("A8").Value = "嘉賓"
$wsContact.Range("B8").Value = "嘉仔"
$wsContact.Range("C8").Value = "9033 2294"

# Column widths (best-fit, approximated through the host's char->pixel rounding)
$wsContact.Columns.Item(1).ColumnWidth = 9.59
$wsContact.Columns.Item(2).ColumnWidth = 8.09
$wsContact.Columns.Item(5).ColumnWidth = 9.25

# Header row is taller (15.75) for row 1
$wsContact.Rows.Item(1).RowHeight = 15.75

$wsContact.Range("D5").Select() | Out-Null

# --- Car sheet data ---
$wsCar.Range("A1").Value = "Self"
$wsCar.Range("B1").Value = "VF 1680"
$wsCar.Range("A2").Value = "Edwin"

$wsCar.Range("B2").Select() | Out-Null

# --- Seating sheet selection change ---
$wsSeating.Range("F15").Select() | Out-Null
